$wb = $excel.ActiveWorkbook

# Helper: set a header/key-style cell (matches existing "key" style used across
# the workbook - black Helvetica Neue 10pt, no border => resolves to existing
# cellXf index 5 instead of minting a new style).
function Set-KeyCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $c.Font.Name = "Helvetica Neue"
    $c.Font.Size = 10
    $c.Font.Color = 0
}

function Set-PlainCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------------
# Info
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Add($null, $lastSheet)
$wsInfo.Name = "Info"

Set-KeyCell   $wsInfo "A1" "apikey"
Set-KeyCell   $wsInfo "B1" "from"
Set-KeyCell   $wsInfo "C1" "max"
Set-KeyCell   $wsInfo "D1" "offset"

Set-KeyCell   $wsInfo "A2" "LaGTXdtMACLYviUe5Is6AB661Lslnded6BmX7eZD"
Set-KeyCell   $wsInfo "B2" 1539146092
Set-PlainCell $wsInfo "C2" 10
Set-PlainCell $wsInfo "D2" 3
[void]$wsInfo.Range("D2").Select()

# ---------------------------------------------------------------------------
# WalletCoinAddrNew
# ---------------------------------------------------------------------------
$wsWCAN = $wb.Worksheets.Add($null, $wsInfo)
$wsWCAN.Name = "WalletCoinAddrNew"

Set-KeyCell   $wsWCAN "A1" "walletID"
Set-KeyCell   $wsWCAN "B1" "info"
Set-PlainCell $wsWCAN "C1" "apikey"

Set-KeyCell   $wsWCAN "A2" "ll1kvbH8C92dwdYlQURF"
Set-PlainCell $wsWCAN "B2" "Testing"
Set-KeyCell   $wsWCAN "C2" "LaGTXdtMACLYviUe5Is6AB661Lslnded6BmX7eZD"
[void]$wsWCAN.Range("C2").Select()

# ---------------------------------------------------------------------------
# WalletCoinNew
# ---------------------------------------------------------------------------
$wsWCN = $wb.Worksheets.Add($null, $wsWCAN)
$wsWCN.Name = "WalletCoinNew"

Set-KeyCell   $wsWCN "A1" "walletName"
Set-PlainCell $wsWCN "B1" "info"
Set-PlainCell $wsWCN "C1" "apikey"

Set-PlainCell $wsWCN "A2" "Coin Wallet"
Set-PlainCell $wsWCN "B2" "Testing"
Set-KeyCell   $wsWCN "C2" "LaGTXdtMACLYviUe5Is6AB661Lslnded6BmX7eZD"
[void]$wsWCN.Range("B2").Select()

# ---------------------------------------------------------------------------
# WalletCoinWithdrawinitiate
# ---------------------------------------------------------------------------
$wsWCWI = $wb.Worksheets.Add($null, $wsWCN)
$wsWCWI.Name = "WalletCoinWithdrawinitiate"

Set-KeyCell   $wsWCWI "A1" "walletID"
Set-KeyCell   $wsWCWI "B1" "amt"
Set-KeyCell   $wsWCWI "C1" "toAddr"
Set-KeyCell   $wsWCWI "D1" "msg"
Set-KeyCell   $wsWCWI "E1" "pin"
Set-KeyCell   $wsWCWI "F1" "apikey"

Set-KeyCell   $wsWCWI "A2" "bDE4PKxaDfblOPc0u6hq"
Set-PlainCell $wsWCWI "B2" 1
Set-KeyCell   $wsWCWI "C2" "1CNugRENyVP6aCbUku6TnryNEs6a41eMKF"
Set-PlainCell $wsWCWI "D2" "Testing"
Set-PlainCell $wsWCWI "E2" 123456
Set-KeyCell   $wsWCWI "F2" "h947NqE3snlyWjznSVFW2UaBLRHzIS62CcY1KhjA"
[void]$wsWCWI.Range("F2").Select()

# ---------------------------------------------------------------------------
# WalletSendToExchange
# ---------------------------------------------------------------------------
$wsWSTE = $wb.Worksheets.Add($null, $wsWCWI)
$wsWSTE.Name = "WalletSendToExchange"

Set-KeyCell   $wsWSTE "A1" "amt"
Set-KeyCell   $wsWSTE "B1" "walletID"
Set-KeyCell   $wsWSTE "C1" "msg"
Set-KeyCell   $wsWSTE "D1" "apikey"

Set-PlainCell $wsWSTE "A2" 1
Set-KeyCell   $wsWSTE "B2" "ll1kvbH8C92dwdYlQURF"
Set-PlainCell $wsWSTE "C2" "Testing"
Set-KeyCell   $wsWSTE "D2" "LaGTXdtMACLYviUe5Is6AB661Lslnded6BmX7eZD"
[void]$wsWSTE.Range("D2").Select()

# ---------------------------------------------------------------------------
# OrderID
# ---------------------------------------------------------------------------
$wsOID = $wb.Worksheets.Add($null, $wsWSTE)
$wsOID.Name = "OrderID"

Set-KeyCell   $wsOID "A1" "apikey"
Set-KeyCell   $wsOID "B1" "orderID"

Set-KeyCell   $wsOID "A2" "LaGTXdtMACLYviUe5Is6AB661Lslnded6BmX7eZD"
Set-PlainCell $wsOID "B2" 12345

# ---------------------------------------------------------------------------
# AskBidNew
# ---------------------------------------------------------------------------
$wsABN = $wb.Worksheets.Add($null, $wsOID)
$wsABN.Name = "AskBidNew"

Set-KeyCell   $wsABN "A1" "rate"
Set-KeyCell   $wsABN "B1" "vol"
Set-KeyCell   $wsABN "C1" "apikey"

Set-KeyCell   $wsABN "A2" 1000000
Set-KeyCell   $wsABN "B2" 1000000
Set-KeyCell   $wsABN "C2" "LaGTXdtMACLYviUe5Is6AB661Lslnded6BmX7eZD"

# This sheet becomes the active/selected tab (activeTab=14, 0-based => 15th sheet).
$wsABN.Activate()
